# تعديل حدث في Card23 - الصف 23
# Row 25's event (A/L/M/N) is cleared - the event it described
# ("هلك حساس سيفتي سويتش باب امامي شمال الكرد" / "تم تغير سيفتي سويتش",
# dated 18/12/2025) is re-homed one row down onto row 26, whose
# previously-blank B:K cells are backfilled with "nan" to match the
# sheet's usual placeholder pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Row 25: blank out the card id and the event/correction/date cells.
$ws.Range("A25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""

# Row 26: fill the previously-empty B:K measurement cells with "nan",
# matching the pattern already used by row 25 (and every other event row).
$ws.Range("B26:K26").Value = "nan"
